# Weekly update: insert a new week's price row for Cilantro (Agrícola del
# Norte S.A. de Arica) at the top of the data block (row 56), pushing the
# existing rows 56-88 down to 57-89. The inserted row reuses the
# commercialisation-unit/origin/min-max-avg price figures that used to sit
# on the old row 56, but carries a new date and a new reported volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 56..88 down to 57..89, duplicating formatting from the row
# that used to be there (matches Excel's default "insert" behaviour).
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with this week's record.
$ws.Cells.Item(56, 1).Value = 1
$ws.Cells.Item(56, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(56, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(56, 4).Value = 44830
$ws.Cells.Item(56, 5).Value = 15
$ws.Cells.Item(56, 6).Value = 100112040
$ws.Cells.Item(56, 7).Value = "Cilantro"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 300
$ws.Cells.Item(56, 11).Value = 1000
$ws.Cells.Item(56, 12).Value = 1200
$ws.Cells.Item(56, 13).Value = 1100
$ws.Cells.Item(56, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(56, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(56, 16).Value = 550
$ws.Cells.Item(56, 17).Value = 2
$ws.Cells.Item(56, 18).Value = "Hortaliza"
